$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the numeric summary fields (items 3-6) and the two free-text
#        fields (7-8) near the top of the sheet ---
$ws.Range("A9").Value  = "3. Số ngày làm việc theo quy định của pháp luật lao động trong tháng: 28"
$ws.Range("A10").Value = "4. Số ngày nghỉ trong tháng (có phép): 2"
$ws.Range("A11").Value = "5. Số ngày nghỉ trong tháng (không phép): 2"
$ws.Range("A12").Value = "6. Số lần vi phạm quy chế, quy định: 1"
$ws.Range("F12").Value = "7. Hành vi vi phạm: aaaa"
$ws.Range("I12").Value = "8. Hình thức kỷ luật: cccc"

# --- 2. Update the "Tên lãnh đạo trực tiếp đánh giá" names for the existing
#        table rows (N16:N18 all shared the same text, N19 its own) ---
$ws.Range("N16").Value = "Bùi Thanh San"
$ws.Range("N17").Value = "Bùi Thanh San"
$ws.Range("N18").Value = "Bùi Thanh San"
$ws.Range("N19").Value = "Nguyễn Đình Long"

# --- 3. Insert a new row 21 so that the existing blank separator row and
#        the two footer rows move down by one, giving us a blank row 20 to
#        fill in as the 5th line item of the detail table ---
$ws.Rows("21").Insert() | Out-Null

# Match the border / alignment formatting used by the rest of the table body
$newRow = $ws.Range("A20:O20")
$newRow.Borders.LineStyle = 1
$newRow.Borders.Color = 0
$newRow.HorizontalAlignment = -4108
$newRow.VerticalAlignment = -4108
$newRow.WrapText = 1

$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Đánh máy"
$ws.Range("C20").Value = "31/03/2025"
$ws.Range("D20").Value = "31/03/2025"
$ws.Range("E20").Value = "05/04/2025"
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = "aaaaa"
$ws.Range("H20").Value = "Vượt tiến độ hoặc có chất lượng hoặc hiệu quả cao"
$ws.Range("I20").Value = "Chưa đánh giá"
$ws.Range("L20").Value = "Chưa phê duyệt"
$ws.Range("M20").Value = 95

# --- 4. Refresh the visible selection so it covers the (now taller) table ---
$ws.Range("A15:O20").Select() | Out-Null
